$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report date range) ---
$ws.Range("A8").Value = "Volume 33   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/2/2026  Through  2/8/2026"

# --- Cells that switch from a number to the "0" / "***.*" placeholder text ---
# (copying an existing placeholder cell brings along both the correct shared-text value and its style)
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("D31").Copy($ws.Range("C31"))

# --- Cells that switch from the "0" / "***.*" placeholder text to a real number ---
$ws.Range("F22").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 3
$ws.Range("H16").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F22").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("H16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 0
$ws.Range("F22").Copy($ws.Range("C33"))
$ws.Range("C33").Value = 1
$ws.Range("F22").Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1
$ws.Range("F22").Copy($ws.Range("I33"))
$ws.Range("I33").Value = 1

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 9
$ws.Range("K16").Value = -30.76923076923
$ws.Range("L16").Value = 28.571428571428
$ws.Range("M16").Value = -60.869565217391
$ws.Range("N16").Value = -90.90909090909
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 16.666666666666
$ws.Range("L17").Value = 110
$ws.Range("M17").Value = 162.5
$ws.Range("N17").Value = 162.5
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -76.190476190476
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = -77.777777777777
$ws.Range("L18").Value = -40
$ws.Range("M18").Value = -53.846153846153
$ws.Range("N18").Value = -96.319018404908
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 43
$ws.Range("K19").Value = 6.976744186046
$ws.Range("L19").Value = 24.324324324324
$ws.Range("M19").Value = -11.538461538461
$ws.Range("N19").Value = -48.888888888888
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -33.333333333333
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = -25
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = -25
$ws.Range("N20").Value = -97.339246119733
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -23.529411764705
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -28.089887640449
$ws.Range("I21").Value = 94
$ws.Range("J21").Value = 119
$ws.Range("K21").Value = -21.008403361344
$ws.Range("L21").Value = 27.027027027027
$ws.Range("M21").Value = -16.071428571428
$ws.Range("N21").Value = -88.437884378843
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -60
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = -36.734693877551
$ws.Range("F24").Value = 143
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = -11.180124223602
$ws.Range("I24").Value = 215
$ws.Range("J24").Value = 213
$ws.Range("K24").Value = 0.93896713615
$ws.Range("L24").Value = 28.742514970059
$ws.Range("M24").Value = 150
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 39
$ws.Range("E25").Value = -43.589743589743
$ws.Range("F25").Value = 110
$ws.Range("G25").Value = 134
$ws.Range("H25").Value = -17.910447761194
$ws.Range("I25").Value = 145
$ws.Range("J25").Value = 177
$ws.Range("K25").Value = -18.079096045197
$ws.Range("L25").Value = 16
$ws.Range("C26").Value = 4
$ws.Range("E26").Value = -63.636363636363
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -18.181818181818
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 37
$ws.Range("K26").Value = -2.702702702702
$ws.Range("L26").Value = 89.473684210526
$ws.Range("M26").Value = 38.461538461538
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("I28").Value = 4
$ws.Range("J28").Value = 4
$ws.Range("L28").Value = 33.333333333333

Write-Output "Edit complete"